# Update the running-time figures from the "triple unrolled init loop"
# measurements and adjust the view selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated measured cycle counts (column B) -- formulas in D/E/F recalc automatically.
$ws.Range("B6").Value = 708
$ws.Range("B8").Value = 355
$ws.Range("B10").Value = 350
$ws.Range("B14").Value = 357

# Update the view: scrolled position and active selection.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D18").Select()
